# Update crypto price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.196.70"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "2.888.78"
$ws.Range("E3").Value = "  -3.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "2.888.34"
$ws.Range("E8").Value = "  -3.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("E11").Value = "  -4.65%  "
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("E13").Value = "  -4.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.66%  "
$ws.Range("D15").Value = "3.388.63"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "60.185.90"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").Value = "2.891.91"
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("E19").Value = "  -4.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "447.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.10%  "
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("E22").Value = "  -6.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.68%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.00%  "
$ws.Range("E32").Value = "  -5.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.12%  "
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0772"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0373"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("D40").Value = "2.857.11"
$ws.Range("E40").Value = "  -10.38%  "
$ws.Range("E41").Value = "  -6.30%  "
$ws.Range("E42").Value = "  -5.04%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -4.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("E48").Value = "  -4.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "111.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.95%  "
$ws.Range("E50").Value = "  -4.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.12%  "
